$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("apiculture")

$ws.Range("C5").Value = "https://raw.githubusercontent.com/AlDenervaud/champdupuits/refs/heads/main/data/images/apiculture/pain_epices.jpg"

$ws.Activate()
$ws.Range("D10").Select()
